{"js": "// M2Doc test fixture update: a Word field \"{ m:self.name }\" is flattened to\n// plain-text runs \"{\", \"m\", \":\", \"self\" (kept orange), \".\", \"name\", \"}\" and a\n// duplicated 4-space run right before \"demonstration\" is dropped.\n\n// ---------------------------------------------------------------------\n// Change 1: remove the extra run of 4 spaces that sits right before the\n// word \"demonstration\" in the first paragraph (there are two identical\n// 4-space runs in the document; we must only touch the second one).\n// ---------------------------------------------------------------------\nconst spaceMatches = context.document.body.search(\"    \", { matchCase: true });\nspaceMatches.load(\"items\");\nawait context.sync();\n\nif (spaceMatches.items.length < 2) {\n  throw new Error(\"Expected at least two runs of four spaces in the document\");\n}\n// The second match (document order) is the one immediately preceding\n// \"demonstration\" and is the one removed by the edit.\nspaceMatches.items[1].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 2: turn the \" m:self.name \" field (fldChar begin/instrText*/\n// fldChar end) into plain <w:t> runs: \"{\", \"m\", \":\", \"self\" (colored),\n// \".\", \"name\", \"}\".\n// ---------------------------------------------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet fieldParagraph = null;\nfor (const p of paragraphs.items) {\n  const fields = p.fields;\n  fields.load(\"items\");\n  await context.sync();\n  if (fields.items.length > 0) {\n    fields.items[0].load(\"code\");\n    await context.sync();\n    if ((fields.items[0].code || \"\").trim() === \"m:self.name\") {\n      fieldParagraph = p;\n      break;\n    }\n  }\n}\n\nif (!fieldParagraph) {\n  throw new Error(\"Could not find the 'm:self.name' field paragraph\");\n}\n\nconst fieldRange = fieldParagraph.getRange();\n\n// Office.js can only inject OOXML using the \"flat OPC\" package wrapper\n// (plain WordprocessingML fragments are rejected), so we build a minimal\n// flat-OPC document containing only the replacement paragraph.\nconst newParagraphXml =\n  '<w:p w:rsidP=\"00F5495F\" w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" ' +\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r w:rsidR=\"00DE6D5A\"><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  '<w:r w:rsidRPr=\"00872F39\">' +\n  '<w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr>' +\n  '<w:t>self</w:t></w:r>' +\n  '<w:r><w:t>.</w:t></w:r>' +\n  '<w:r w:rsidR=\"00F0712B\"><w:t>name</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships></pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  newParagraphXml +\n  '</w:body></w:document></pkg:xmlData></pkg:part>' +\n  '</pkg:package>';\n\nfieldRange.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# M2Doc test fixture update: a Word field \"{ m:self.name }\" is flattened to\n# plain-text runs \"{\", \"m\", \":\", \"self\" (kept orange), \".\", \"name\", \"}\" and a\n# duplicated 4-space run right before \"demonstration\" is dropped.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: remove the extra run of 4 spaces that sits right before the\n# word \"demonstration\" in the first paragraph (there are two identical\n# 4-space runs in the document; only the second one - right after\n# \"NotExistingEPackage\" - is removed).\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"    \"\n$find.Forward = $true\n$find.Wrap = 0          # wdFindStop - do not wrap back to start\n\n$matchIndex = 0\n$targetRange = $null\nwhile ($find.Execute()) {\n    $matchIndex = $matchIndex + 1\n    if ($matchIndex -eq 2) {\n        $targetRange = $d.Content.Duplicate\n        $targetRange.Start = $find.Parent.Start\n        $targetRange.End = $find.Parent.End\n        break\n    }\n    $find.Parent.Collapse(0)   # wdCollapseEnd, keep searching forward\n}\n\nif ($matchIndex -lt 2) {\n    throw \"Expected at least two runs of four spaces in the document\"\n}\n\n$find.Parent.Delete()\n\n# ---------------------------------------------------------------------\n# Change 2: turn the \" m:self.name \" field (fldChar begin/instrText*/\n# fldChar end) into plain <w:t> runs: \"{\", \"m\", \":\", \"self\" (colored),\n# \".\", \"name\", \"}\".\n# ---------------------------------------------------------------------\n$field = $null\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n    $candidate = $d.Fields.Item($i)\n    if ($candidate.Code.Text.Trim() -eq \"m:self.name\") {\n        $field = $candidate\n        break\n    }\n}\n\nif ($field -eq $null) {\n    throw \"Could not find the 'm:self.name' field\"\n}\n\n# Locate the paragraph that contains the field so we can replace its whole\n# content (fldChar begin .. fldChar end) in one shot.\n$fieldStart = $field.Code.Start\n$fieldParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($fieldStart -ge $p.Range.Start -and $fieldStart -lt $p.Range.End) {\n        $fieldParagraph = $p\n        break\n    }\n}\n\nif ($fieldParagraph -eq $null) {\n    throw \"Could not locate the paragraph containing the field\"\n}\n\n$newParagraphXml = '<w:p w:rsidP=\"00F5495F\" w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>{</w:t></w:r><w:r w:rsidR=\"00DE6D5A\"><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r w:rsidRPr=\"00872F39\"><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>self</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidR=\"00F0712B\"><w:t>name</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p>'\n\n$fieldParagraph.Range.InsertXML($newParagraphXml)\n"}
